# Apply daily crypto price/symbol-list refresh (GitHub Actions update, 2022-12-19).
# Updates Price (column D) figures for several coins and realigns the
# Coin/Link/Volume columns (B, C, E) for rows whose ranking position changed,
# matching the canonical OOXML diff for this commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "241.86"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.158"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.05595"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "6.395"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.8046"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9550"
$ws.Range("B10").Value = "One"
$ws.Range("C10").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0005752"
$ws.Range("E10").Value = "9OneONE"
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1407"
$ws.Range("E11").Value = "10WazirXWRX"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07374"
$ws.Range("E12").Value = "11MandalaExchangeTokenMDX"
$ws.Range("B13").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C13").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03100"
$ws.Range("E13").Value = "12LiechtensteinCryptoassetsExchangeLCX"
$ws.Range("B14").Value = "BitrueCoin"
$ws.Range("C14").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.03056"
$ws.Range("E14").Value = "13BitrueCoinBTR"
$ws.Range("B15").Value = "BitMartToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.09279"
$ws.Range("E15").Value = "14BitMartTokenBMX"
$ws.Range("B16").Value = "MCDex"
$ws.Range("C16").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.627"
$ws.Range("E16").Value = "15MCDexMCB"
$ws.Range("B17").Value = "BitForexToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.001651"
$ws.Range("E17").Value = "16BitForexTokenBF"
$ws.Range("B18").Value = "CoinExToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.04693"
$ws.Range("E18").Value = "17CoinExTokenCET"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.006372"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.004982"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0001500"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.756"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.3260"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.03876"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006866"
$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.003401"
$ws.Range("E42").Value = "41CEJICEJI"
$ws.Range("B43").Value = "BKEXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1034"
$ws.Range("E43").Value = "42BKEXTokenBKK"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.008403"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005935"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0005502"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.6826"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.09997"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002100"